$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44911, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Bing", "Primera", 650, 500, 600, 554, "`$/kilo (en caja de 15 kilos)", "Región del Maule", 554, 1),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44911, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Bing", "Segunda", 250, 400, 400, 400, "`$/kilo (en caja de 15 kilos)", "Región del Maule", 400, 1),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44911, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Primera", 450, 500, 600, 556, "`$/kilo (en caja de 15 kilos)", "Región del Maule", 556, 1),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44911, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Segunda", 185, 400, 400, 400, "`$/kilo (en caja de 15 kilos)", "Región del Maule", 400, 1),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44911, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Rainier", "Primera", 450, 700, 800, 744, "`$/kilo (en caja de 15 kilos)", "Región del Maule", 744, 1),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44911, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Rainier", "Segunda", 150, 500, 500, 500, "`$/kilo (en caja de 15 kilos)", "Región del Maule", 500, 1),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44911, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Santina", "Primera", 200, 600, 600, 600, "`$/kilo (en caja de 15 kilos)", "Región del Maule", 600, 1),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44911, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Santina", "Segunda", 155, 400, 400, 400, "`$/kilo (en caja de 15 kilos)", "Región del Maule", 400, 1),
)

$startRow = 285
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}